$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '26.421.64'
$ws.Range('E2').Value = '  -1.10%  '
$ws.Range('D3').Value = '1.624.75'
$ws.Range('E3').Value = '  -0.68%  '
$ws.Range('D4').Value = '1.01'
$ws.Range('E4').Value = '  +0.35%  '
$ws.Range('D5').Value = '212.66'
$ws.Range('E5').Value = '  -0.57%  '
$ws.Range('D6').Value = '0.495'
$ws.Range('E6').Value = '  +0.70%  '
$ws.Range('D7').Value = '1.01'
$ws.Range('E7').Value = '  +0.39%  '
$ws.Range('D8').Value = '0.249'
$ws.Range('E8').Value = '  -1.23%  '
$ws.Range('D9').Value = '0.0620'
$ws.Range('E9').Value = '  -0.18%  '
$ws.Range('D10').Value = '18.92'
$ws.Range('E10').Value = '  -0.90%  '
$ws.Range('D11').Value = '0.0838'
$ws.Range('E11').Value = '  +0.39%  '
$ws.Range('D12').Value = '1.851.27'
$ws.Range('E12').Value = '  -0.75%  '
$ws.Range('D13').Value = '1.647.50'
$ws.Range('E13').Value = '  +0.73%  '
$ws.Range('D14').Value = '4.10'
$ws.Range('E14').Value = '  +0.97%  '
$ws.Range('D15').Value = '0.520'
$ws.Range('E15').Value = '  -0.99%  '
$ws.Range('D16').Value = '63.99'
$ws.Range('E16').Value = '  +1.05%  '
$ws.Range('D17').Value = '26.412.47'
$ws.Range('E17').Value = '  -1.08%  '
$ws.Range('D18').Value = '0.0₃0738'
$ws.Range('E18').Value = '  +0.50%  '
$ws.Range('D19').Value = '215.20'
$ws.Range('E19').Value = '  +2.73%  '
$ws.Range('D20').Value = '1.01'
$ws.Range('E20').Value = '  +0.42%  '
$ws.Range('D21').Value = '4.29'
$ws.Range('E21').Value = '  -1.00%  '
$ws.Range('D22').Value = '6.21'
$ws.Range('E22').Value = '  +1.78%  '
$ws.Range('D23').Value = '9.27'
$ws.Range('E23').Value = '  -1.63%  '
$ws.Range('D24').Value = '1.98'
$ws.Range('E24').Value = '  +4.78%  '
$ws.Range('D25').Value = '147.82'
$ws.Range('E25').Value = '  +1.36%  '
$ws.Range('D26').Value = '1.01'
$ws.Range('E26').Value = '  +0.38%  '
$ws.Range('D27').Value = '0.120'
$ws.Range('E27').Value = '  -0.93%  '
$ws.Range('D28').Value = '6.82'
$ws.Range('E28').Value = '  +2.34%  '
$ws.Range('D29').Value = '15.54'
$ws.Range('E29').Value = '  +0.83%  '
$ws.Range('D30').Value = '0.0506'
$ws.Range('E30').Value = '  -2.49%  '
$ws.Range('D31').Value = '1.16'
$ws.Range('E31').Value = '  -1.58%  '
$ws.Range('D32').Value = '3.31'
$ws.Range('E32').Value = '  +2.56%  '
$ws.Range('D33').Value = '2.93'
$ws.Range('E33').Value = '  -0.82%  '
$ws.Range('D34').Value = '1.49'
$ws.Range('E34').Value = '  -1.30%  '
$ws.Range('D35').Value = '2.38'
$ws.Range('E35').Value = '  -1.50%  '
$ws.Range('D36').Value = '1.212.74'
$ws.Range('E36').Value = '  +3.89%  '
$ws.Range('D37').Value = '0.0173'
$ws.Range('E37').Value = '  +2.87%  '
$ws.Range('D38').Value = '1.01'
$ws.Range('E38').Value = '  +0.45%  '
$ws.Range('D39').Value = '0.794'
$ws.Range('E39').Value = '  -2.46%  '
$ws.Range('D40').Value = '0.499'
$ws.Range('E40').Value = '  -0.95%  '
$ws.Range('D41').Value = '2.25'
$ws.Range('E41').Value = '  -3.35%  '
$ws.Range('D42').Value = '0.792'
$ws.Range('E42').Value = '  +0.32%  '
$ws.Range('D43').Value = '5.36'
$ws.Range('E43').Value = '  -1.05%  '
$ws.Range('D44').Value = '1.760.04'
$ws.Range('E44').Value = '  -0.68%  '
$ws.Range('D45').Value = '92.49'
$ws.Range('E45').Value = '  -0.06%  '
$ws.Range('D46').Value = '1.57'
$ws.Range('E46').Value = '  +1.11%  '
$ws.Range('D47').Value = '54.66'
$ws.Range('E47').Value = '  -0.05%  '
$ws.Range('D48').Value = '0.0₆0102'
$ws.Range('E48').Value = '  -0.81%  '
$ws.Range('D49').Value = '0.0511'
$ws.Range('E49').Value = '  -0.16%  '
$ws.Range('D50').Value = '7.61'
$ws.Range('E50').Value = '  -0.56%  '
$ws.Range('D51').Value = '0.408'
$ws.Range('E51').Value = '  -0.59%  '
